$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($addr, $val) {
    $rng = $ws.Range($addr)
    $rng.NumberFormat = "@"
    $rng.Value = $val
}

Set-TextValue "E2"  "-1.03%"

Set-TextValue "D3"  "27.18"
Set-TextValue "E3"  "3.22%"

Set-TextValue "D4"  "5.101"
Set-TextValue "E4"  "0.15%"

Set-TextValue "D5"  "0.05666"
Set-TextValue "E5"  "0.92%"

Set-TextValue "E6"  "-0.48%"

Set-TextValue "D7"  "0.8231"
Set-TextValue "E7"  "1.33%"

Set-TextValue "D8"  "0.8467"
Set-TextValue "E8"  "0.41%"

Set-TextValue "D9"  "0.1327"
Set-TextValue "E9"  "-1.35%"

Set-TextValue "D10" "0.06924"
Set-TextValue "E10" "-0.45%"

Set-TextValue "D11" "0.02885"
Set-TextValue "E11" "2.29%"

Set-TextValue "E12" "-0.07%"

Set-TextValue "D13" "0.001516"
Set-TextValue "E13" "0.31%"

Set-TextValue "D14" "0.04121"
Set-TextValue "E14" "-12.13%"

Set-TextValue "D15" "0.0005979"
Set-TextValue "E15" "-0.26%"

Set-TextValue "D16" "0.006120"
Set-TextValue "E16" "0.11%"

Set-TextValue "E17" "-1.57%"

Set-TextValue "D18" "3.001"
Set-TextValue "E18" "-1.73%"

Set-TextValue "D19" "2.226"
Set-TextValue "E19" "5.09%"

Set-TextValue "D21" "0.03137"
Set-TextValue "E21" "-0.07%"

Set-TextValue "D22" "0.1291"
Set-TextValue "E22" "-0.71%"

Set-TextValue "D23" "3.558"
Set-TextValue "E23" "-5.51%"

Set-TextValue "E24" "-0.01%"

Set-TextValue "E25" "-2.52%"

Set-TextValue "D26" "0.004458"
Set-TextValue "E26" "-3.48%"

Set-TextValue "D28" "0.0001438"
Set-TextValue "E28" "3.55%"

Set-TextValue "E40" "0.20%"

Set-TextValue "D41" "0.006039"
Set-TextValue "E41" "-0.92%"

Set-TextValue "D42" "0.1054"
Set-TextValue "E42" "-0.20%"

Set-TextValue "D43" "0.002510"
Set-TextValue "E43" "0.48%"

Set-TextValue "D44" "0.008298"
Set-TextValue "E44" "-4.46%"

Set-TextValue "D45" "0.00005319"
Set-TextValue "E45" "0.53%"

Set-TextValue "E46" "0.07%"

Set-TextValue "D47" "0.1010"
Set-TextValue "E47" "-36.83%"

Set-TextValue "D48" "0.002589"
Set-TextValue "E48" "25.46%"

Set-TextValue "E49" "0.07%"

Set-TextValue "E50" "0.07%"
